$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Table 4.1")
$ws2 = $wb.Worksheets.Item("Table 4.2")

# ---------------------------------------------------------------------
# Table 4.1 ("ws1") - update Guessed_Distance (D) and Team 1/2 measurement
# (F/G) columns for rows 3-12, and fill in the 11th-object guess (D13).
# The H column (MEDIAN(E:G)) is a formula and recalculates automatically.
# ---------------------------------------------------------------------
$ws1.Range("D3").Value2 = 45
$ws1.Range("D4").Value2 = 28
$ws1.Range("D5").Value2 = 12
$ws1.Range("D6").Value2 = 125
$ws1.Range("D7").Value2 = 175
$ws1.Range("D8").Value2 = 500
$ws1.Range("D9").Value2 = 6
$ws1.Range("D10").Value2 = 35
$ws1.Range("D11").Value2 = 146
$ws1.Range("D12").Value2 = 255
$ws1.Range("D13").Value2 = 131

$ws1.Range("F3").Value2 = 41.7
$ws1.Range("F4").Value2 = 27.9
$ws1.Range("F5").Value2 = 15.11
$ws1.Range("F6").Value2 = 82.7
$ws1.Range("F7").Value2 = 138
$ws1.Range("F8").Value2 = 273
$ws1.Range("F9").Value2 = 8
$ws1.Range("F10").Value2 = 46
$ws1.Range("F11").Value2 = 106.7
$ws1.Range("F12").Value2 = 197.7

$ws1.Range("G3").Value2 = 39
$ws1.Range("G4").Value2 = 27
$ws1.Range("G5").Value2 = 16
$ws1.Range("G6").Value2 = 84
$ws1.Range("G7").Value2 = 130.8
$ws1.Range("G8").Value2 = 260.4
$ws1.Range("G9").Value2 = 7.1
$ws1.Range("G10").Value2 = 46
$ws1.Range("G11").Value2 = 106.7
$ws1.Range("G12").Value2 = 91

# ---------------------------------------------------------------------
# Table 4.2 ("ws2") - Linear model block (rows 4-13)
# D4:D13 is the Guessed_Distance input column; it previously held the
# integer "0" number format (style used nowhere else on this sheet) -
# paste the General-format style already used by the equivalent column
# on Table 4.1 before writing the new numbers, matching the workbook
# author re-typing the values with the default format.
# ---------------------------------------------------------------------
$ws1.Range("D3:D12").Copy()
$ws2.Range("D4").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws2.Range("D4").Value2 = 45
$ws2.Range("D5").Value2 = 28
$ws2.Range("D6").Value2 = 17
$ws2.Range("D7").Value2 = 83
$ws2.Range("D8").Value2 = 110
$ws2.Range("D9").Value2 = 330
$ws2.Range("D10").Value2 = 3
$ws2.Range("D11").Value2 = 59
$ws2.Range("D12").Value2 = 120
$ws2.Range("D13").Value2 = 265

# 11th-object placeholder value in the linear block is cleared out.
$ws2.Range("D14").ClearContents()

# Regression coefficients (Solver output) for the linear fit.
$ws2.Range("L4").Value2 = 9.812
$ws2.Range("N4").Value2 = 0.7778

# O4 re-typed (same visible text, "*Gussed_Distance").
$ws2.Range("O4").Value2 = "*Gussed_Distance"

# ---------------------------------------------------------------------
# Table 4.2 - Quadratic model block (rows 21-30): only the Solver
# coefficients changed; the Guessed_Distance inputs (D21:D30) stayed put.
# ---------------------------------------------------------------------
$ws2.Range("L21").Value2 = 3.888299
$ws2.Range("N21").Value2 = 0.9269883
$ws2.Range("Q21").Value2 = -0.0004515

# ---------------------------------------------------------------------
# Workbook-level defined name used by Solver to track the optimized cell
# for the quadratic model: it now points one row down (H33 instead of
# H32, i.e. the "Error Standard deviation" cell instead of the SSE sum).
# ---------------------------------------------------------------------
$solverOpt = $wb.Names.Item("Table 4.2!solver_opt")
$solverOpt.RefersTo = "='Table 4.2'!`$H`$33"

# ---------------------------------------------------------------------
# Restore sheet selections/scroll position to match the saved view.
# ---------------------------------------------------------------------
$ws2.Activate()
$ws2.Range("F47").Select()

$ws1.Activate()
$ws1.Range("E3:G3").Select()
